$d = $word.ActiveDocument
$wns = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

function Set-ParaXml($para, $innerXml) {
    $xml = '<w:p xmlns:w="' + $wns + '">' + $innerXml + '</w:p>'
    $para.Range.InsertXML($xml)
}

function Add-ParaXml($innerXml) {
    $p = $d.Paragraphs.Add()
    Set-ParaXml $p $innerXml
}

# --- 1) Paragraph 3 "Enrutamiento y puertos:" - drop <w:lang w:val="en-US"/> ---
Set-ParaXml $d.Paragraphs.Item(3) '<w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Enrutamiento y puertos</w:t></w:r><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>:</w:t></w:r>'

# --- 2) "Consume menos anchoo de banda" -> "Consume menos ancho de banda" ---
$d.Content.Find.Execute("anchoo de banda", $true, $false, $false, $false, $false, $true, 1, $false, "ancho de banda", 2) | Out-Null

# --- 3) Drop trailing curly quote after "...en la red." ---
$d.Content.Find.Execute([string][char]0x2019, $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# --- 4) "viaja pero no el contenido" -> split runs + add comma ---
$pAntes = $d.Paragraphs.Item(46)
Set-ParaXml $pAntes '<w:r><w:t xml:space="preserve">Antes del 2018 la mayor&#237;a del tr&#225;fico no viajaba encriptado, pero en la actualidad la mayor&#237;a de nuestros paquetes de datos viajan protegidos y solo es visible a d&#243;nde </w:t></w:r><w:r><w:t>viaja,</w:t></w:r><w:r><w:t xml:space="preserve"> pero no el contenido.</w:t></w:r>'

# --- 5) "reciben todos lo paquetes" -> split runs + fix "lo" -> "los" ---
$pISP = $d.Paragraphs.Item(47)
Set-ParaXml $pISP '<w:r><w:t xml:space="preserve">Los </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>ISP</w:t></w:r><w:r><w:t xml:space="preserve"> reciben todos </w:t></w:r><w:r><w:t xml:space="preserve">los </w:t></w:r><w:r><w:t>paquetes de datos que enviamos y los env&#237;an a su destino, pudiendo aplicar filtros de normativas referentes a bloquear p&#225;ginas de pirater&#237;a, contenido protegido geogr&#225;ficamente o pol&#237;ticas gubernamentales seg&#250;n lo que determine cada gobierno.</w:t></w:r>'

# --- 6) "publicas" gets wrapped in proofErr spell-check markers ---
$pVPN = $d.Paragraphs.Item(49)
Set-ParaXml $pVPN '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:t>Otro m&#233;todo es una</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> VPN</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">(red privada virtual) </w:t></w:r><w:r><w:t xml:space="preserve">tecnolog&#237;a que permite una extensi&#243;n segura de una red local sobre una red </w:t></w:r><w:r><w:t>p&#250;blica</w:t></w:r><w:r><w:t xml:space="preserve"> como Internet, permite que nuestra computadora </w:t></w:r><w:r><w:t>env&#237;e</w:t></w:r><w:r><w:t xml:space="preserve"> y reciba dato conect&#225;ndose a otras redes compartidas o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>publicas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> como si fueran una red privada con </w:t></w:r><w:r><w:t>toda la funcionalidad</w:t></w:r><w:r><w:t xml:space="preserve">, seguridad, y </w:t></w:r><w:r><w:t>pol&#237;ticas de gesti&#243;n de una red privada</w:t></w:r><w:r><w:t>.</w:t></w:r>'

# --- 7) Append the new "examen grupo 6" presentation content at the end of the doc ---

# Two blank bold paragraphs
Add-ParaXml '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>'
Add-ParaXml '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>'

# Q: El enrutamiento o routing es:
Add-ParaXml '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">El enrutamiento o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>routing</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> es:</w:t></w:r>'

# A: Elegir la mejor ruta...
Add-ParaXml '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:t>Elegir la mejor ruta para enviar datos desde una red a otra</w:t></w:r>'

# Q: La conexion entre dos computadoras...
Add-ParaXml '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>La conexi&#243;n entre dos computadoras depende de m&#250;ltiples factores, por lo tanto, nunca se conectan a trav&#233;s de las mismas direcciones IP.</w:t></w:r>'

# A: Falso
Add-ParaXml '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:t>Falso</w:t></w:r>'

# Q: Que utiliza nuestra computadora...
Add-ParaXml '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>&#191;Qu&#233; utiliza nuestra computadora para diferenciar solicitudes a diferentes aplicaciones que se encuentran en un mismo servidor?</w:t></w:r>'

# A: Puertos
Add-ParaXml '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:t>Puertos</w:t></w:r>'

# Q: Que es el ISP?
Add-ParaXml '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>&#191;Qu&#233; es el ISP?</w:t></w:r>'

# A: El proveedor de servicios de Internet
Add-ParaXml '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:t>El proveedor de servicios de Internet</w:t></w:r>'

# Q: Las tablas de enrutamiento...
Add-ParaXml '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Las tablas de enrutamiento nos sirven para determinar qu&#233; camino deben seguir los datos</w:t></w:r>'

# A: Verdadero
Add-ParaXml '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:t>Verdadero</w:t></w:r>'

# Q: El Proxy y las VPN son PROTOCOLOS...
Add-ParaXml '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>El Proxy y las VPN son PROTOCOLOS que se utilizan para EL DIRECCIONAMIENTO</w:t></w:r>'

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
